# Daily attendance processing - reorders the list of "Changed By" users
# in column G (Session Analysis Results sheet) by reversing the order of
# the comma-separated entries in each cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 7).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 157 }

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value()
    if ($null -ne $val -and $val -is [string] -and $val.Contains(",")) {
        $parts = $val -split ","
        $trimmed = @()
        foreach ($p in $parts) { $trimmed += $p.Trim() }

        $reversed = @()
        for ($i = $trimmed.Count - 1; $i -ge 0; $i--) {
            $reversed += $trimmed[$i]
        }

        $cell.Value = [string]::Join(", ", $reversed)
    }
}
